$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.082.65"
$ws.Range("E2").Value = "'  -0.19%  "
$ws.Range("D3").Value = "'1.637.35"
$ws.Range("E3").Value = "'  -1.88%  "
$ws.Range("E4").Value = "'  -0.15%  "
$ws.Range("D5").Value = "'214.02"
$ws.Range("E5").Value = "'  +1.64%  "
$ws.Range("D6").Value = "'0.5258"
$ws.Range("E6").Value = "'  +0.10%  "
$ws.Range("D8").Value = "'0.2598"
$ws.Range("E8").Value = "'  -1.33%  "
$ws.Range("D9").Value = "'0.06312"
$ws.Range("E9").Value = "'  +0.20%  "
$ws.Range("D10").Value = "'20.70"
$ws.Range("E10").Value = "'  -2.35%  "
$ws.Range("E11").Value = "'  +1.22%  "
$ws.Range("D12").Value = "'1.627.49"
$ws.Range("E12").Value = "'  -2.51%  "
$ws.Range("D13").Value = "'4.423"
$ws.Range("E13").Value = "'  -0.48%  "
$ws.Range("D14").Value = "'1.861.46"
$ws.Range("E14").Value = "'  -1.95%  "
$ws.Range("D15").Value = "'0.5499"
$ws.Range("E15").Value = "'  -1.30%  "
$ws.Range("D16").Value = "'0.0₅8169"
$ws.Range("E16").Value = "'  +3.05%  "
$ws.Range("D17").Value = "'65.06"
$ws.Range("E17").Value = "'  -2.75%  "
$ws.Range("D18").Value = "'26.073.56"
$ws.Range("E19").Value = "'  -0.09%  "
$ws.Range("D20").Value = "'4.692"
$ws.Range("E20").Value = "'  -1.15%  "
$ws.Range("D21").Value = "'188.14"
$ws.Range("E21").Value = "'  +0.82%  "
$ws.Range("D22").Value = "'10.14"
$ws.Range("E22").Value = "'  -2.31%  "
$ws.Range("D23").Value = "'6.163"
$ws.Range("E23").Value = "'  -0.31%  "
$ws.Range("E24").Value = "'  -0.20%  "
$ws.Range("E25").Value = "'  -2.41%  "
$ws.Range("E26").Value = "'  -2.76%  "
$ws.Range("D27").Value = "'7.416"
$ws.Range("E27").Value = "'  -1.31%  "
$ws.Range("D28").Value = "'15.85"
$ws.Range("E28").Value = "'  -0.71%  "
$ws.Range("D29").Value = "'1.411"
$ws.Range("E29").Value = "'  +4.42%  "
$ws.Range("D30").Value = "'0.06024"
$ws.Range("E30").Value = "'  -3.69%  "
$ws.Range("E31").Value = "'  -1.93%  "
$ws.Range("D32").Value = "'3.444"
$ws.Range("E32").Value = "'  -2.18%  "
$ws.Range("D33").Value = "'3.411"
$ws.Range("E33").Value = "'  -0.30%  "
$ws.Range("D34").Value = "'1.639"
$ws.Range("E34").Value = "'  +0.45%  "
$ws.Range("D35").Value = "'0.9868"
$ws.Range("E35").Value = "'  -1.32%  "
$ws.Range("E36").Value = "'  +1.29%  "
$ws.Range("D37").Value = "'2.396"
$ws.Range("E37").Value = "'  -0.64%  "
$ws.Range("E38").Value = "'  -5.13%  "
$ws.Range("D39").Value = "'0.01619"
$ws.Range("E39").Value = "'  +0.17%  "
$ws.Range("D40").Value = "'0.8566"
$ws.Range("E40").Value = "'  -1.96%  "
$ws.Range("B41").Value = "'PaxDollar"
$ws.Range("C41").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.001"
$ws.Range("E41").Value = "'  -0.24%  "
$ws.Range("B42").Value = "'FraxShare"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.729"
$ws.Range("E42").Value = "'  -6.70%  "
$ws.Range("B43").Value = "'Maker"
$ws.Range("C43").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'1.038.74"
$ws.Range("E43").Value = "'  -5.69%  "
$ws.Range("D44").Value = "'100.67"
$ws.Range("E44").Value = "'  +0.35%  "
$ws.Range("D45").Value = "'1.787.49"
$ws.Range("E45").Value = "'  -1.90%  "
$ws.Range("D46").Value = "'0.0₈108"
$ws.Range("E46").Value = "'  -1.94%  "
$ws.Range("D47").Value = "'55.61"
$ws.Range("E47").Value = "'  +0.52%  "
$ws.Range("D48").Value = "'1.003"
$ws.Range("E48").Value = "'  -0.20%  "
$ws.Range("D49").Value = "'8.049"
$ws.Range("E49").Value = "'  +0.13%  "
$ws.Range("D50").Value = "'0.05172"
$ws.Range("E50").Value = "'  -1.22%  "
$ws.Range("D51").Value = "'0.4222"
$ws.Range("E51").Value = "'  -0.60%  "
